$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.184.68'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.077.08'
$ws.Range('E3').Value = '  -0.69%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.94'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +10.53%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.393'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '61.49'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.36%  '
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '16.31'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.73%  '
$ws.Range('D14').Value = '2.379.51'
$ws.Range('E14').Value = '  -0.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.820'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.83%  '
$ws.Range('D17').Value = '2.076.41'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '37.231.90'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.45'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +13.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '74.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').Value = '0.0₃0934'
$ws.Range('E21').Value = '  +10.35%  '
$ws.Range('E22').Value = '  +5.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.79'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +14.77%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.41'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.48'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('E30').Value = '  +3.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.14'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.80'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0640'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.74%  '
$ws.Range('E34').Value = '  +9.41%  '
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.33'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.09%  '
$ws.Range('B37').Value = 'BinanceUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +27.10%  '
$ws.Range('E39').Value = '  -3.86%  '
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0228'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.39'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('E46').Value = '  +1.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.62'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +14.39%  '
$ws.Range('E48').Value = '  +7.30%  '
$ws.Range('D49').Value = '1.311.06'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.98'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.16%  '
